# Generate Report for Handback
# The localization handback for 3537563b-77e7-4663-aa2e-957d039a8b2b.md has come in,
# so it swaps places with f5b1a7b0-c673-41c6-a019-b923d6cabc77.md in every report
# table (rows are sorted by file name), and its status / timestamps are refreshed.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet "Overview"
# ---------------------------------------------------------------------------
$ov = $wb.Worksheets.Item("Overview")

$ov.Range("A2").Value = '3537563b-77e7-4663-aa2e-957d039a8b2b.md'
$ov.Range("B2").Value = 'e2e\3537563b-77e7-4663-aa2e-957d039a8b2b.md'
$ov.Range("G2").Value = '2016-08-13 23:03:56'

$ov.Range("A3").Value = 'f5b1a7b0-c673-41c6-a019-b923d6cabc77.md'
$ov.Range("B3").Value = 'e2e\f5b1a7b0-c673-41c6-a019-b923d6cabc77.md'
$ov.Range("E3").Value = 'Handed back: in sync with en-US'
$ov.Range("F3").Value = 'Handed back: in sync with en-US'
$ov.Range("G3").Value = '2016-08-13 23:02:46'

$ov.Hyperlinks.Delete()
$ov.Hyperlinks.Add($ov.Range("B2"), "https://github.com/OpenLocalizationTestOrg/oltest/blob/9920891c7827685675df55308e6ea0986f86b542/e2e/f5b1a7b0-c673-41c6-a019-b923d6cabc77.md", "", "", 'e2e\3537563b-77e7-4663-aa2e-957d039a8b2b.md')
$ov.Hyperlinks.Add($ov.Range("B3"), "https://github.com/OpenLocalizationTestOrg/oltest/blob/29169a5af2eae7852a2f2554e367736faaab700b/e2e/3537563b-77e7-4663-aa2e-957d039a8b2b.md", "", "", 'e2e\f5b1a7b0-c673-41c6-a019-b923d6cabc77.md')

# ---------------------------------------------------------------------------
# Sheet "zh-cn"
# ---------------------------------------------------------------------------
$zh = $wb.Worksheets.Item("zh-cn")

$zh.Range("A2").Value = '3537563b-77e7-4663-aa2e-957d039a8b2b.md'
$zh.Range("G2").Value = '3537563b-77e7-4663-aa2e-957d039a8b2b.85932e4da4ec79813dd07a19961c04b08e0771c4.zh-cn.xlf'
$zh.Range("H2").Value = '2016-08-13 23:03:49'
$zh.Range("I2").Value = '3537563b-77e7-4663-aa2e-957d039a8b2b.md'
$zh.Range("J2").Value = '3537563b-77e7-4663-aa2e-957d039a8b2b.85932e4da4ec79813dd07a19961c04b08e0771c4.zh-cn.xlf'
$zh.Range("K2").Value = '2016-08-13 23:04:17'

$zh.Range("A3").Value = 'f5b1a7b0-c673-41c6-a019-b923d6cabc77.md'
$zh.Range("C3").Value = 'Handed back: in sync with en-US'
$zh.Range("G3").Value = 'f5b1a7b0-c673-41c6-a019-b923d6cabc77.5dbadb4022f65fb9cf506737aa36012f487f99e8.zh-cn.xlf'
$zh.Range("H3").Value = '2016-08-13 23:02:38'
$zh.Range("I3").Value = 'f5b1a7b0-c673-41c6-a019-b923d6cabc77.md'
$zh.Range("J3").Value = 'f5b1a7b0-c673-41c6-a019-b923d6cabc77.5dbadb4022f65fb9cf506737aa36012f487f99e8.zh-cn.xlf'
$zh.Range("P3").Value = ''

$zh.Hyperlinks.Delete()
$zh.Hyperlinks.Add($zh.Range("A2"), "https://github.com/OpenLocalizationTestOrg/oltest/blob/9920891c7827685675df55308e6ea0986f86b542/e2e/f5b1a7b0-c673-41c6-a019-b923d6cabc77.md", "", "", '3537563b-77e7-4663-aa2e-957d039a8b2b.md')
$zh.Hyperlinks.Add($zh.Range("I2"), "https://github.com/OpenLocalizationTestOrg/ol-test0-zhcn/blob/9395db14ebe5d8a35217ee31a503782ae7b9a6eb/e2e/f5b1a7b0-c673-41c6-a019-b923d6cabc77.md", "", "", '3537563b-77e7-4663-aa2e-957d039a8b2b.md')
$zh.Hyperlinks.Add($zh.Range("A3"), "https://github.com/OpenLocalizationTestOrg/oltest/blob/29169a5af2eae7852a2f2554e367736faaab700b/e2e/3537563b-77e7-4663-aa2e-957d039a8b2b.md", "", "", 'f5b1a7b0-c673-41c6-a019-b923d6cabc77.md')
$zh.Hyperlinks.Add($zh.Range("I3"), "https://github.com/OpenLocalizationTestOrg/ol-test0-zhcn/blob/9395db14ebe5d8a35217ee31a503782ae7b9a6eb/e2e/3537563b-77e7-4663-aa2e-957d039a8b2b.md", "", "", 'f5b1a7b0-c673-41c6-a019-b923d6cabc77.md')

$zh.Columns.Item(16).ColumnWidth = 13.7470528738839

# ---------------------------------------------------------------------------
# Sheet "de-de"
# ---------------------------------------------------------------------------
$de = $wb.Worksheets.Item("de-de")

$de.Range("A2").Value = '3537563b-77e7-4663-aa2e-957d039a8b2b.md'
$de.Range("G2").Value = '3537563b-77e7-4663-aa2e-957d039a8b2b.85932e4da4ec79813dd07a19961c04b08e0771c4.de-de.xlf'
$de.Range("H2").Value = '2016-08-13 23:03:56'
$de.Range("I2").Value = '3537563b-77e7-4663-aa2e-957d039a8b2b.md'
$de.Range("J2").Value = '3537563b-77e7-4663-aa2e-957d039a8b2b.85932e4da4ec79813dd07a19961c04b08e0771c4.de-de.xlf'
$de.Range("K2").Value = '2016-08-13 23:04:27'

$de.Range("A3").Value = 'f5b1a7b0-c673-41c6-a019-b923d6cabc77.md'
$de.Range("C3").Value = 'Handed back: in sync with en-US'
$de.Range("G3").Value = 'f5b1a7b0-c673-41c6-a019-b923d6cabc77.5dbadb4022f65fb9cf506737aa36012f487f99e8.de-de.xlf'
$de.Range("H3").Value = '2016-08-13 23:02:46'
$de.Range("I3").Value = 'f5b1a7b0-c673-41c6-a019-b923d6cabc77.md'
$de.Range("J3").Value = 'f5b1a7b0-c673-41c6-a019-b923d6cabc77.5dbadb4022f65fb9cf506737aa36012f487f99e8.de-de.xlf'
$de.Range("P3").Value = ''

$de.Hyperlinks.Delete()
$de.Hyperlinks.Add($de.Range("A2"), "https://github.com/OpenLocalizationTestOrg/oltest/blob/9920891c7827685675df55308e6ea0986f86b542/e2e/f5b1a7b0-c673-41c6-a019-b923d6cabc77.md", "", "", '3537563b-77e7-4663-aa2e-957d039a8b2b.md')
$de.Hyperlinks.Add($de.Range("I2"), "https://github.com/OpenLocalizationTestOrg/ol-test0-dede/blob/ae2b9082b4fb2633f62a9edb46278ae367622957/e2e/f5b1a7b0-c673-41c6-a019-b923d6cabc77.md", "", "", '3537563b-77e7-4663-aa2e-957d039a8b2b.md')
$de.Hyperlinks.Add($de.Range("A3"), "https://github.com/OpenLocalizationTestOrg/oltest/blob/29169a5af2eae7852a2f2554e367736faaab700b/e2e/3537563b-77e7-4663-aa2e-957d039a8b2b.md", "", "", 'f5b1a7b0-c673-41c6-a019-b923d6cabc77.md')
$de.Hyperlinks.Add($de.Range("I3"), "https://github.com/OpenLocalizationTestOrg/ol-test0-dede/blob/ae2b9082b4fb2633f62a9edb46278ae367622957/e2e/3537563b-77e7-4663-aa2e-957d039a8b2b.md", "", "", 'f5b1a7b0-c673-41c6-a019-b923d6cabc77.md')

$de.Columns.Item(16).ColumnWidth = 13.7470528738839
